$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Append three new food-item rows (289-291) to the existing table.
$ws.Cells.Item(289, 1).Value = "caramel popcorn with cashews"
$ws.Cells.Item(289, 2).Value = "0.33 cup"
$ws.Cells.Item(289, 3).Value = 120
$ws.Cells.Item(289, 4).Value = 1
$ws.Cells.Item(289, 5).Value = 6
$ws.Cells.Item(289, 6).Value = 16

$ws.Cells.Item(290, 1).Value = "skinny popcorn"
$ws.Cells.Item(290, 2).Value = "1 cup"
$ws.Cells.Item(290, 3).Value = 40
$ws.Cells.Item(290, 4).Value = 0.5329999999999999
$ws.Cells.Item(290, 5).Value = 2.67
$ws.Cells.Item(290, 6).Value = 4

$ws.Cells.Item(291, 1).Value = "kala chana onion tomato salad"
$ws.Cells.Item(291, 2).Value = "1 serving"
$ws.Cells.Item(291, 3).Value = 75.845
$ws.Cells.Item(291, 4).Value = 5.657500000000001
$ws.Cells.Item(291, 5).Value = 0.99
$ws.Cells.Item(291, 6).Value = 10.7075
